# First block of moving DumpShort functions
# Re-colour the SI unit table by category, and add new :DEFAULT/:ALL/:electrical/
# :mechanical/:special marker rows, then update the "litre" footnote text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- colour palette (fgColor of a "solid" Interior pattern) ----
# theme 9 (0-based) / tint 0.8  -> Excel ThemeColor index 10
# theme 8 (0-based) / tint 0.8  -> Excel ThemeColor index 9
# rgb FFFCFFA8 (yellow)          -> Color 11075580
# rgb FFE4B5DB (pink)            -> Color 14398948

$themeGreen = 10
$themeBlue  = 9
$yellow = 11075580
$pink   = 14398948

function Set-Cell($addr, $value, $themeColorIdx, $rgbColor, $bold) {
    $c = $ws.Range($addr)
    $c.Value = $value
    if ($bold) {
        $c.Font.Name = "Menlo"
    }
    if ($themeColorIdx -ne $null) {
        $c.Interior.ThemeColor = $themeColorIdx
    } elseif ($rgbColor -ne $null) {
        $c.Interior.Color = $rgbColor
    }
}

# Row 3 - metre / second base units -> green
$ws.Range("B3").Interior.ThemeColor = $themeGreen
$ws.Range("C3").Interior.Color = $pink
$ws.Range("B4").Interior.ThemeColor = $themeGreen
$ws.Range("C4").Interior.Color = $pink
$ws.Range("B5").Interior.ThemeColor = $themeGreen
$ws.Range("C5").Interior.Color = $pink
$ws.Range("B6").Interior.Color = $yellow
$ws.Range("C6").Interior.Color = $pink
$ws.Range("B7").Interior.Color = $pink
$ws.Range("C7").Interior.Color = $pink
$ws.Range("B8").Interior.Color = $pink
$ws.Range("C8").Interior.Color = $yellow
$ws.Range("B9").Interior.ThemeColor = $themeBlue
$ws.Range("C9").Interior.Color = $yellow
$ws.Range("C10").Interior.Color = $yellow
$ws.Range("C11").Interior.Color = $yellow
$ws.Range("C12").Interior.Color = $yellow
$ws.Range("C13").Interior.Color = $yellow

$ws.Range("B14").Value = ":DEFAULT"
$ws.Range("B14").Interior.ThemeColor = $themeGreen
$ws.Range("C14").Interior.Color = $yellow

$ws.Range("B15").Value = ":electrical"
$ws.Range("B15").Interior.Color = $yellow
$ws.Range("C15").Interior.Color = $yellow

$ws.Range("B16").Value = ":mechanical"
$ws.Range("B16").Interior.Color = $pink
$ws.Range("C16").Interior.ThemeColor = $themeBlue

$ws.Range("B17").Value = ":special"
$ws.Range("B17").Interior.ThemeColor = $themeBlue
$ws.Range("C17").Interior.ThemeColor = $themeBlue

$ws.Range("B18").Value = ":ALL"
$ws.Range("C18").Interior.ThemeColor = $themeBlue

$ws.Range("C19").Interior.ThemeColor = $themeBlue
$ws.Range("C20").Interior.ThemeColor = $themeBlue
$ws.Range("C21").Interior.ThemeColor = $themeBlue
$ws.Range("C22").Interior.ThemeColor = $themeGreen

$ws.Range("C23").Value = "#litre added due to common use of ml, dl, etc. "

$ws.Range("B20").Select()
